# Generate Report for Archive
#
# 1. Change the status text "Ready for handoff" -> "In Translation"
#    wherever it appears (Overview!E2/F2, zh-cn!C2, de-de!C2).
# 2. Narrow the "Status" columns (Overview columns E & F, and column C
#    on the zh-cn / de-de sheets) from their old width to the new,
#    narrower width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text ------------------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Narrow the status columns ------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
